$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6194122874597942
$ws.Range("C3").Value = 0.6188812961859822
$ws.Range("D3").Value = 3
$ws.Range("C4").Value = 0.6205972201845776
$ws.Range("C5").Value = 0.6122253332066087
